# The experiment's image-path strings were edited to reflect the "huang"
# test folder on Windows: forward slashes became backslashes, and the
# second filename changed from "CS+4" to "CS-3".
#
# Column A rows 2-11 all shared the string "Extinction/CS+3.BMP" and rows
# 12-21 all shared "Extinction/CS+4.BMP". Re-assigning each cell's Value
# updates (and de-duplicates) the underlying shared-string table so every
# cell in a block keeps pointing at a single shared string, exactly like
# the original layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "Extinction\CS+3.BMP"
}

for ($r = 12; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = "Extinction\CS-3.BMP"
}

# The sheet's selection moved from D5 to A12:A21 (active cell A12).
$ws.Range("A12:A21").Select()
